$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (Price / Volume(1h)) values.
# Each target cell is forced to Text format before assignment so the
# numeric-looking / percent-looking strings are stored verbatim (matching
# the original inline-string cells) instead of being auto-converted to
# numbers by Excel's smart entry. The style is then reset to Normal so
# no stray cell formatting is left behind.

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "300.84"
Set-TextValue "E2" "-0.01%"
Set-TextValue "D3" "32.35"
Set-TextValue "E3" "2.00%"
Set-TextValue "D4" "4.964"
Set-TextValue "E4" "-2.58%"
Set-TextValue "D5" "0.07664"
Set-TextValue "E5" "-2.11%"
Set-TextValue "D6" "1.920"
Set-TextValue "E6" "-17.35%"
Set-TextValue "D7" "7.836"
Set-TextValue "E7" "0.43%"
Set-TextValue "D8" "3.797"
Set-TextValue "E8" "-0.78%"
Set-TextValue "D9" "0.9178"
Set-TextValue "E9" "0.35%"
Set-TextValue "D10" "0.1748"
Set-TextValue "E10" "-0.63%"
Set-TextValue "D11" "0.07767"
Set-TextValue "E11" "2.33%"
Set-TextValue "D12" "0.08580"
Set-TextValue "E12" "-6.68%"
Set-TextValue "D13" "0.03163"
Set-TextValue "E13" "2.47%"
Set-TextValue "D14" "0.09986"
Set-TextValue "E14" "-0.36%"
Set-TextValue "E15" "-0.65%"
Set-TextValue "D16" "0.005932"
Set-TextValue "E16" "2.22%"
Set-TextValue "D17" "3.468"
Set-TextValue "E17" "-0.58%"
Set-TextValue "D18" "2.153"
Set-TextValue "E18" "-4.14%"
Set-TextValue "E19" "1.89%"
Set-TextValue "E20" "-0.88%"
Set-TextValue "D21" "4.290"
Set-TextValue "E21" "6.40%"
Set-TextValue "D22" "0.1991"
Set-TextValue "E22" "11.17%"
Set-TextValue "D23" "0.04511"
Set-TextValue "E23" "-1.85%"
Set-TextValue "D24" "0.001222"
Set-TextValue "E24" "-2.36%"
Set-TextValue "D25" "0.004401"
Set-TextValue "E25" "-1.38%"
Set-TextValue "D26" "0.0001251"
Set-TextValue "E26" "0.00%"
Set-TextValue "E39" "-3.78%"
Set-TextValue "D40" "0.04692"
Set-TextValue "E40" "-2.69%"
Set-TextValue "D41" "0.007482"
Set-TextValue "E41" "3.17%"
Set-TextValue "D42" "0.1349"
Set-TextValue "E42" "-0.81%"
Set-TextValue "D43" "0.002332"
Set-TextValue "E43" "6.48%"
Set-TextValue "D44" "0.01050"
Set-TextValue "E44" "2.18%"
Set-TextValue "D45" "0.00006244"
Set-TextValue "E45" "0.94%"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "E46" "0.05%"
Set-TextValue "E47" "10.10%"
Set-TextValue "D49" "0.00002101"
Set-TextValue "E49" "0.05%"
Set-TextValue "D50" "0.0002001"
Set-TextValue "E50" "0.05%"
